$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: blank spacer row, formatted the same as row 16 (copy format down) ---
$ws.Range("A16:D16").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)

# --- Row 18: new KPI atomic-name-fix entry ---
# Values are written left-to-right (A, B, C, D) so new shared-string entries
# get minted in the same order as the source workbook.
$ws.Range("A18").Value = "PoS 2019 - MT Conv Big - REG"

# B18 reuses the existing formatted look from a similar cell elsewhere in the
# sheet (B11), then gets the 0.00 number format applied on top - matching the
# new style introduced for this row.
$ws.Range("B11").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B18").Value = "Ice Tea Shelf: Top Shelf"
$ws.Range("B18").NumberFormat = "0.00"

$ws.Range("C18").Value = "Ice Tea Shelf: Fuze Berry-Hibiscus- 1L"

# D18 reuses the existing formatted look from D2, then gets the 0.00 number
# format applied on top - matching the new style introduced for this row.
$ws.Range("D2").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = "Ice Tea Shelf: Fuze Berry-Hibiscus - 1L"
$ws.Range("D18").NumberFormat = "0.00"

$ws.Range("E18").Formula = "=CONCATENATE(""UPDATE ``static``.atomic_kpi a JOIN ``static``.kpi k ON k.pk=a.kpi_fk JOIN ``static``.kpi_set s ON s.pk=k.kpi_set_fk   SET a.name='"",D18,""', a.description='"",D18,""', a.display_text='"",D18,""'  WHERE s.name='"",A18,""' AND k.display_text='"",B18,""' AND a.name='"",C18,""';"")"

# The wrap-text style used by C16/C17 round-trips through this engine without
# its wrapText flag (source file was authored in a different tool), so make
# sure the wrap formatting is still explicitly applied to C17. Done last so
# it doesn't shift the cellXf indices minted above for B18/D18.
$ws.Range("C17").WrapText = $true

# Match the author's final selection (Excel leaves the last-used cell selected)
$ws.Range("E18").Select()
